$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.581.24'
$ws.Range("E2").Value = '  -2.38%  '

$ws.Range("D3").Value = '2.295.96'
$ws.Range("E3").Value = '  -4.72%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '545.61'
$ws.Range("E5").Value = '  -1.59%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.90'
$ws.Range("E6").Value = '  -3.91%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("E8").Value = '  -3.00%  '

$ws.Range("D9").Value = '2.295.10'
$ws.Range("E9").Value = '  -4.72%  '

$ws.Range("E10").Value = '  -2.88%  '

$ws.Range("E11").Value = '  -2.11%  '

$ws.Range("E12").Value = '  +1.14%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.333'
$ws.Range("E13").Value = '  -5.43%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.88'
$ws.Range("E14").Value = '  -3.44%  '

$ws.Range("D15").Value = '2.704.02'
$ws.Range("E15").Value = '  -4.84%  '

$ws.Range("D16").Value = '58.586.28'
$ws.Range("E16").Value = '  -2.18%  '

$ws.Range("E17").Value = '  -3.54%  '

$ws.Range("D18").Value = '2.239.95'
$ws.Range("E18").Value = '  -6.98%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.64'
$ws.Range("E19").Value = '  -4.86%  '

$ws.Range("E20").Value = '  -4.67%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '313.69'
$ws.Range("E21").Value = '  -3.87%  '

$ws.Range("E22").Value = '  -4.45%  '

$ws.Range("E23").Value = '  +0.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.28'
$ws.Range("E24").Value = '  -2.25%  '

$ws.Range("E25").Value = '  -6.04%  '

$ws.Range("E26").Value = '  -0.14%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.08'
$ws.Range("E27").Value = '  -6.11%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.31'
$ws.Range("E28").Value = '  -5.80%  '

$ws.Range("E29").Value = '  -2.58%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '170.46'
$ws.Range("E30").Value = '  -0.40%  '

$ws.Range("D31").Value = '0.0₃0723'
$ws.Range("E31").Value = '  -6.03%  '

$ws.Range("E32").Value = '  -0.17%  '

$ws.Range("E33").Value = '  -5.58%  '

$ws.Range("E34").Value = '  -5.63%  '

$ws.Range("E35").Value = '  -0.03%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.75'
$ws.Range("E36").Value = '  -3.50%  '

$ws.Range("E37").Value = '  -0.01%  '

$ws.Range("E38").Value = '  -5.79%  '

$ws.Range("E39").Value = '  -6.20%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '38.06'
$ws.Range("E40").Value = '  -2.27%  '

$ws.Range("E41").Value = '  -5.49%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '297.67'
$ws.Range("E42").Value = '  -7.67%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '140.05'
$ws.Range("E43").Value = '  -4.60%  '

$ws.Range("E44").Value = '  -4.89%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0949'
$ws.Range("E45").Value = '  -1.38%  '

$ws.Range("E46").Value = '  -3.13%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.553'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.39'
$ws.Range("E48").Value = '  -7.18%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0214'
$ws.Range("E49").Value = '  -3.06%  '

$ws.Range("E50").Value = '  -4.31%  '

$ws.Range("E51").Value = '  -0.28%  '
